# Weekly fruit/vegetable update: a new daily price record for
# "Poroto granado" at Vega Central Mapocho de Santiago was inserted
# as row 356, pushing the existing rows 356:385 down to 357:386.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 356 (shifts 356:385 -> 357:386).
$ws.Rows.Item(356).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(356, 1).Value = 9
$ws.Cells.Item(356, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(356, 3).Value = "Metropolitana"
$ws.Cells.Item(356, 4).Value = 45021
$ws.Cells.Item(356, 5).Value = 13
$ws.Cells.Item(356, 6).Value = 100112030
$ws.Cells.Item(356, 7).Value = "Poroto granado"
$ws.Cells.Item(356, 8).Value = "Sin especificar"
$ws.Cells.Item(356, 9).Value = "Primera"
$ws.Cells.Item(356, 10).Value = 75
$ws.Cells.Item(356, 11).Value = 32000
$ws.Cells.Item(356, 12).Value = 34000
$ws.Cells.Item(356, 13).Value = 33200
$ws.Cells.Item(356, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(356, 15).Value = "Región Metropolitana"
$ws.Cells.Item(356, 16).Value = 1328
$ws.Cells.Item(356, 17).Value = 25
$ws.Cells.Item(356, 18).Value = "Hortaliza"
